$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42 (shifting the existing rows 42-83 down to 43-84)
$ws.Rows.Item(42).Insert(-4121, $null)

# Fill in the new row 42 with this week's price report
# (same market/product/region as the surrounding rows, new date + updated prices)
$ws.Cells.Item(42, 1).Value2 = 10
$ws.Cells.Item(42, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(42, 3).Value2 = "La Araucanía"
$ws.Cells.Item(42, 4).Value2 = 44966
$ws.Cells.Item(42, 5).Value2 = 9
$ws.Cells.Item(42, 6).Value2 = "Fruta"
$ws.Cells.Item(42, 7).Value2 = 100108
$ws.Cells.Item(42, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(42, 9).Value2 = 100108004
$ws.Cells.Item(42, 10).Value2 = "Papaya"
$ws.Cells.Item(42, 11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(42, 12).Value2 = "Primera"
$ws.Cells.Item(42, 13).Value2 = 50
$ws.Cells.Item(42, 14).Value2 = 30000
$ws.Cells.Item(42, 15).Value2 = 30000
$ws.Cells.Item(42, 16).Value2 = 30000
$ws.Cells.Item(42, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(42, 18).Value2 = "Provincia del Elquí"
$ws.Cells.Item(42, 19).Value2 = 3000
$ws.Cells.Item(42, 20).Value2 = 10
